# styled dataframe display mode
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B column values for existing rows 2-8 (scaled display values)
$ws.Range("B2").Value = -55.01340103149414
$ws.Range("B3").Value = -59.02640151977539
$ws.Range("B4").Value = -43.45640182495117
$ws.Range("B5").Value = 49.46620178222656
$ws.Range("B6").Value = 86.60639953613281
$ws.Range("B7").Value = 65.66690063476562
$ws.Range("B8").Value = 98.10030364990234

# Append a new summary row 9 ("Пастбище")
$ws.Range("A9").Value = "Пастбище"
$ws.Range("B9").Value = 142.3435974121094
$ws.Range("C9").Value = 0.0182
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.3282999992370605
$ws.Range("F9").Value = -0.05119999870657921
